$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Gradient Boosting (DeepWalk)
$ws.Range("B3").Value = 0.8063704052585373
$ws.Range("C3").Value = 0.1662004957965591
$ws.Range("D3").Value = 0.8139999999999999
$ws.Range("E3").Value = 2074738269.309195
$ws.Range("F3").Value = 0.009928829375639398

# Row 5: Linear Regression (DeepWalk)
$ws.Range("B5").Value = 0.6424436849792784
$ws.Range("C5").Value = 0.3201680102268146
$ws.Range("D5").Value = 0.574
$ws.Range("E5").Value = 2819354430.266713
$ws.Range("F5").Value = 0.03448963229494235

# Row 7: Random Forest (DeepWalk)
$ws.Range("B7").Value = 0.8233898036780043
$ws.Range("C7").Value = 0.1640556571329358
$ws.Range("D7").Value = 0.802
$ws.Range("E7").Value = 1981460121.101597
$ws.Range("F7").Value = 0.00955756944128665
